$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the values of columns D, M, N, O, P, S between row 2 <-> row 6
# and row 3 <-> row 7 (weekly Fruta/Hortaliza price records swapped places).

$cols = @("D", "M", "N", "O", "P", "S")
$rowPairs = @(@(2, 6), @(3, 7))

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
